# IPA.xlsx - "new test cases aree added to IPA module"
#
# Adds 8 new Technology/Company "IPA Module" test-case rows (IPA115, IPA02..IPA08)
# to the bottom of the "Test Cases" sheet (rows 31-38), right after the
# existing data block (rows 1-30).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3 already carries the plain bordered style (no fill, no wrap) that the
# new rows use for every column -- clone its formatting onto the new block so
# the new cells pick up the same cell style already present in the workbook
# instead of fabricating a brand-new one.
$ws.Range("A3:E3").Copy()
$ws.Range("A31:E38").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# TCID, Jira id, Description, Runmode ("Y") for each new test case.
# Column E (Results) is intentionally left blank, matching the other
# not-yet-run rows.

# Row 31 - IPA115 placeholder entry
$ws.Range("A31").Value2 = "IPA115"
$ws.Range("B31").Value2 = "OBT"
$ws.Range("C31").Value2 = "OBT"
$ws.Range("D31").Value2 = "Y"

# Row 32 - IPA02: Technology competitors visualization
$ws.Range("A32").Value2 = "IPA02"
$ws.Range("B32").Value2 = "OPQA-4402||OPQA-4403||OPQA-4404||OPQA-4405"
$ws.Range("C32").Value2 = "Veify Technology competitors visualization"
$ws.Range("D32").Value2 = "Y"

# Row 33 - IPA03: Technology trending visualization
$ws.Range("A33").Value2 = "IPA03"
$ws.Range("B33").Value2 = "OPQA-4397||OPQA-4398||OPQA-4400||OPQA-4401"
$ws.Range("C33").Value2 = "Verify Technology trending visualization "
$ws.Range("D33").Value2 = "Y"

# Row 34 - IPA04: Company Technology trending visualization
$ws.Range("A34").Value2 = "IPA04"
$ws.Range("B34").Value2 = "OPQA-4412||OPQA-4413||OPQA-4415"
$ws.Range("C34").Value2 = "Verify Company Technology trending visualization "
$ws.Range("D34").Value2 = "Y"

# Row 35 - IPA05: Company Key Information
$ws.Range("A35").Value2 = "IPA05"
$ws.Range("B35").Value2 = "OPQA-4425||OPQA-4423||OPQA-4421"
$ws.Range("C35").Value2 = "Verify Company Key Information"
$ws.Range("D35").Value2 = "Y"

# Row 36 - IPA06: Technology Key Information
$ws.Range("A36").Value2 = "IPA06"
$ws.Range("B36").Value2 = "OPQA-4444||OPQA-4445||OPQA-4446"
$ws.Range("C36").Value2 = "Verify Technology Key Information"
$ws.Range("D36").Value2 = "Y"

# Row 37 - IPA07: Mandatory Field and Sorting - Technology ResultList
$ws.Range("A37").Value2 = "IPA07"
$ws.Range("B37").Value2 = "OPQA-4387||OPQA-4372||OPQA-4373||OPQA-4374||OPQA-4376||OPQA-4377||OPQA-4378||OPQA-4379"
$ws.Range("C37").Value2 = "Verify Mandatory Field and Sorting with Different options at Technology ResultList"
$ws.Range("D37").Value2 = "Y"

# Row 38 - IPA08: Mandatory Field and Sorting - Company ResultList
$ws.Range("A38").Value2 = "IPA08"
$ws.Range("B38").Value2 = "OPQA-4387||OPQA-4380||OPQA-4381||OPQA-4384||OPQA-4385||OPQA-4386"
$ws.Range("C38").Value2 = "Verify Mandatory Field and Sorting with Different options at Company ResultList"
$ws.Range("D38").Value2 = "Y"

# Select the whole table (now A1:E38) so the saved view reflects the newly
# extended used range, the same way the editor's selection grew to cover the
# freshly-added rows.
$ws.Range("A1:E38").Select()
